# Allow turning RoR component visualisation on/off in script C
# Adds a new "calibration_only" parameter row to the "General parameters" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General parameters")

# New row 9: calibration_only parameter
$ws.Range("A9").Value = "calibration_only"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = "select this option to not run hydro-solar-wind scenarios but only calibration run under regular hydropower operation (0 = no, 1 = yes)"

# Match formatting of the other parameter rows (B: Input style, C: wrapped Note style)
$ws.Range("B9").Style = "Input"
$ws.Range("C9").WrapText = $true

# Widen column C to fit the new, slightly longer description text
$ws.Columns.Item(3).ColumnWidth = 60.8

# Row heights reflow because of the wider column / new wrapped row
$ws.Rows.Item(6).RowHeight = 58
$ws.Rows.Item(9).RowHeight = 29

# Cursor ended up on cell E6 on this sheet after the edit
$ws.Range("E6").Select() | Out-Null

# Cursor also moved on the "Hydropower plant parameters" sheet
$ws1 = $wb.Worksheets.Item("Hydropower plant parameters")
$ws1.Range("A9").Select() | Out-Null
